$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the changed "Coin" / "Link" / "Price" / "Volume(1h)" cells to the
# latest scraped snapshot. Price cells that render as plain decimal strings
# (e.g. "318.30") are forced to Text number format first so COM does not
# silently coerce them to numbers and drop the significant trailing zero,
# matching how the source data stores every column as text.
$ws.Range('D2').Value = '44.046.65'
$ws.Range('E2').Value = '  +1.62%  '
$ws.Range('D3').Value = '2.246.46'
$ws.Range('E3').Value = '  +0.40%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '318.30'
$ws.Range('E5').Value = '  -0.07%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '100.98'
$ws.Range('E6').Value = '  +1.92%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.572'
$ws.Range('E7').Value = '  -1.64%  '
$ws.Range('E8').Value = '  +0.12%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.543'
$ws.Range('E9').Value = '  -3.41%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.81'
$ws.Range('E10').Value = '  -0.61%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0825'
$ws.Range('E11').Value = '  -0.63%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.51'
$ws.Range('E12').Value = '  -2.43%  '
$ws.Range('E13').Value = '  -1.93%  '
$ws.Range('D14').Value = '2.591.52'
$ws.Range('E14').Value = '  +0.55%  '
$ws.Range('B15').Value = 'Polygon'
$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.846'
$ws.Range('E15').Value = '  -2.35%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '2.251.25'
$ws.Range('E16').Value = '  +0.31%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.16'
$ws.Range('E17').Value = '  -1.45%  '
$ws.Range('D18').Value = '43.998.42'
$ws.Range('E18').Value = '  +1.58%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.29'
$ws.Range('E19').Value = '  -6.40%  '
$ws.Range('D20').Value = '0.0₃0974'
$ws.Range('E20').Value = '  +0.20%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.42'
$ws.Range('E21').Value = '  -3.44%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '65.49'
$ws.Range('E22').Value = '  +0.31%  '
$ws.Range('E23').Value = '  -3.08%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '234.78'
$ws.Range('E24').Value = '  -0.72%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.06'
$ws.Range('E25').Value = '  -5.45%  '
$ws.Range('E26').Value = '  +0.02%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.49'
$ws.Range('E27').Value = '  +4.39%  '
$ws.Range('E28').Value = '  -0.28%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '37.79'
$ws.Range('E29').Value = '  +2.50%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.09'
$ws.Range('E30').Value = '  -4.63%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '159.26'
$ws.Range('E31').Value = '  +0.95%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '20.04'
$ws.Range('E32').Value = '  -1.15%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0843'
$ws.Range('E33').Value = '  -3.35%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.68'
$ws.Range('E34').Value = '  -1.07%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.18'
$ws.Range('E36').Value = '  +7.46%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.93'
$ws.Range('E37').Value = '  +1.86%  '
$ws.Range('E38').Value = '  -2.31%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '16.04'
$ws.Range('E39').Value = '  +11.34%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.67'
$ws.Range('E40').Value = '  -0.72%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.14'
$ws.Range('E41').Value = '  -5.85%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0314'
$ws.Range('E42').Value = '  -2.35%  '
$ws.Range('E43').Value = '  +0.07%  '
$ws.Range('D44').Value = '1.740.06'
$ws.Range('E44').Value = '  -4.34%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.196'
$ws.Range('E45').Value = '  -3.07%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '81.81'
$ws.Range('E46').Value = '  -2.64%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '74.28'
$ws.Range('E47').Value = '  +0.58%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '5.15'
$ws.Range('E48').Value = '  -3.07%  '
$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '102.36'
$ws.Range('E49').Value = '  -1.02%  '
$ws.Range('B50').Value = 'Stacks'
$ws.Range('C50').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.66'
$ws.Range('E50').Value = '  +3.57%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '57.45'
$ws.Range('E51').Value = '  -1.07%  '
